$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")
$ws.Range("C1").Value = "allow_choice_duplicates"
$ws.Range("C1").WrapText = $true
$ws.Range("C2").Value = "Yes"
$ws.Select()
$ws.Range("C2").Select()
